$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change regenerates the handback-status report: the three tracked
# files (Ping.md, test-scenario-1.md, osmond-test-1.md) are now listed with
# Ping.md first, and Ping.md's handback timestamp/status were refreshed
# (it is no longer in sync with en-US because a newer handback arrived).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "Ping.md"
$wsOverview.Range("B2").Value = "Handed back: not in sync with en-US"

$wsOverview.Range("A3").Value = "test-scenario-1.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A4").Value = "osmond-test-1.md"
$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"

# Hyperlinks keep pointing at the same targets they always did (rId2/3/4),
# only the row they decorate - and hence their display text - moves.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a5287759ee19d999486c30f70b7686d01ea0d695/e2e/test-scenario-1.md", "", "", "Ping.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f4b3096e48dfeb999c383c55a46706a9a3b95935/e2e/osmond-test-1.md", "", "", "test-scenario-1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3e270db8fe1a93168c64ab5fd0c0a7f513156d82/e2e/Ping.md", "", "", "osmond-test-1.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2 -> Ping.md (status flipped to "not in sync", handback datetime refreshed)
$wsZhCn.Range("A2").Value = "Ping.md"
$wsZhCn.Range("B2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("C2").Value = "Ping.95b58797b320fac8f901ea501ac186551882a36c.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-02-01 07:05:43"
$wsZhCn.Range("E2").Value = "Ping.md"
$wsZhCn.Range("F2").Value = "Ping.95b58797b320fac8f901ea501ac186551882a36c.zh-cn.xlf"
$wsZhCn.Range("G2").Value = "2016-02-22 08:37:26"
$wsZhCn.Range("H2").Value = "Include"

# Row 3 -> test-scenario-1.md
$wsZhCn.Range("A3").Value = "test-scenario-1.md"
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "test-scenario-1.974d721459da0ff5eab675ae57b2ea10d235c32c.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2015-11-06 07:23:55"
$wsZhCn.Range("E3").Value = "test-scenario-1.md"
$wsZhCn.Range("F3").Value = "test-scenario-1.974d721459da0ff5eab675ae57b2ea10d235c32c.de-de.xlf"
$wsZhCn.Range("G3").Value = "2015-11-06 07:36:12"
$wsZhCn.Range("H3").Value = "Include"

# Row 4 -> osmond-test-1.md
$wsZhCn.Range("A4").Value = "osmond-test-1.md"
$wsZhCn.Range("B4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C4").Value = "osmond-test-1.f11b754a5dfc36a6a88309daf550f7387cc99a9b.zh-cn.xlf"
$wsZhCn.Range("D4").Value = "2015-11-06 07:23:55"
$wsZhCn.Range("E4").Value = "osmond-test-1.md"
$wsZhCn.Range("F4").Value = "osmond-test-1.f11b754a5dfc36a6a88309daf550f7387cc99a9b.de-de.xlf"
$wsZhCn.Range("G4").Value = "2015-11-06 07:36:12"
$wsZhCn.Range("H4").Value = "Include"

# Hyperlinks: same underlying targets as before (rId2..rId13 in the same
# order), only re-labelled/re-ordered to follow the new row order.
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a5287759ee19d999486c30f70b7686d01ea0d695/e2e/test-scenario-1.md", "", "", "Ping.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9409763d56a0eab0f0ee07cc0c6e2ec61787719d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/master/test-scenario-1.974d721459da0ff5eab675ae57b2ea10d235c32c.zh-cn.xlf", "", "", "Ping.95b58797b320fac8f901ea501ac186551882a36c.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ed2810913f07afe80cd8049db04fd7ec7d55a1c2/e2e/test-scenario-1.md", "", "", "Ping.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e13eb945b0aaf2756c9b50d122f74bb1ed6120a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/master/test-scenario-1.974d721459da0ff5eab675ae57b2ea10d235c32c.de-de.xlf", "", "", "Ping.95b58797b320fac8f901ea501ac186551882a36c.zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f4b3096e48dfeb999c383c55a46706a9a3b95935/e2e/osmond-test-1.md", "", "", "test-scenario-1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9409763d56a0eab0f0ee07cc0c6e2ec61787719d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/master/osmond-test-1.f11b754a5dfc36a6a88309daf550f7387cc99a9b.zh-cn.xlf", "", "", "test-scenario-1.974d721459da0ff5eab675ae57b2ea10d235c32c.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ed2810913f07afe80cd8049db04fd7ec7d55a1c2/e2e/osmond-test-1.md", "", "", "test-scenario-1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e13eb945b0aaf2756c9b50d122f74bb1ed6120a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/master/osmond-test-1.f11b754a5dfc36a6a88309daf550f7387cc99a9b.de-de.xlf", "", "", "test-scenario-1.974d721459da0ff5eab675ae57b2ea10d235c32c.de-de.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3e270db8fe1a93168c64ab5fd0c0a7f513156d82/e2e/Ping.md", "", "", "osmond-test-1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16a291f0e64d4dcc27577d91cf96df62a4779888/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/master/Ping.95b58797b320fac8f901ea501ac186551882a36c.zh-cn.xlf", "", "", "osmond-test-1.f11b754a5dfc36a6a88309daf550f7387cc99a9b.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e90f9d5f500b5cdfb41138847b002b3cb4c280e0/e2e/Ping.md", "", "", "osmond-test-1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/06dba2df2e2eb2b371299de18336731207efc15b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/master/Ping.95b58797b320fac8f901ea501ac186551882a36c.zh-cn.xlf", "", "", "osmond-test-1.f11b754a5dfc36a6a88309daf550f7387cc99a9b.de-de.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de" - header-only sheet, no data rows changed.
# ---------------------------------------------------------------------------
